# Commit: "Sheet names are updated."
# The only semantic change in the diff is the worksheet's (tab) name,
# from the Turkish default "Sayfa1" to "Conversations".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Conversations"
